# Data Siswa.xlsx - Update Besar - Increment 2
# Remove the "Nomor Telepon" (Phone Number) column from the student sheet.
# Before: A=NISN | B=Nama Siswa | C=Nomor Telepon | D=Status
# After:  A=NISN | B=Nama Siswa | C=Status
#
# Deleting column C shifts the old "Status" column (D) left into C and
# drops the phone-number data/header entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Delete()
